# Auto-generated edit script: apply scheduled-runner market data updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 100
$ws.Range("I31").Value = 100
$ws.Range("K31").Value = 300
$ws.Range("M31").Value = -70
$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = 10000
$ws.Range("N40").Value = -10350
$ws.Range("H41").Value = 490
$ws.Range("I41").Value = 336.25
$ws.Range("J41").Value = 900
$ws.Range("K41").Value = 336.25
$ws.Range("L41").Value = 900
$ws.Range("M41").Value = 103.75
$ws.Range("N41").Value = -1780
$ws.Range("H51").Value = 7415.6665
$ws.Range("I51").Value = 7415.6665
$ws.Range("K51").Value = 7415.6665
$ws.Range("M51").Value = -6931.6665
$ws.Range("H58").Value = 825
$ws.Range("I58").Value = 350
$ws.Range("J58").Value = 3200
$ws.Range("K58").Value = 1050
$ws.Range("L58").Value = 9600
$ws.Range("M58").Value = -900
$ws.Range("N58").Value = -9900
$ws.Range("H61").Value = 505
$ws.Range("I61").Value = 505
$ws.Range("K61").Value = 1515
$ws.Range("M61").Value = -1343
$ws.Range("H64").Value = 12500
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 20000
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 20000
$ws.Range("M64").Value = -4752
$ws.Range("N64").Value = -20496
$ws.Range("H67").Value = 12500
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 20000
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 20000
$ws.Range("M67").Value = -4142
$ws.Range("N67").Value = -21716
$ws.Range("H70").Value = 1880.6875
$ws.Range("I70").Value = 1399.5
$ws.Range("K70").Value = 4198.5
$ws.Range("M70").Value = -3928.5
$ws.Range("H73").Value = 1880.6875
$ws.Range("I73").Value = 1399.5
$ws.Range("K73").Value = 4198.5
$ws.Range("M73").Value = -3262.5
$ws.Range("H74").Value = 5237.25
$ws.Range("I74").Value = 5237.25
$ws.Range("K74").Value = 5237.25
$ws.Range("M74").Value = -4301.25
$ws.Range("H76").Value = 3828.1428
$ws.Range("I76").Value = 3632.8333
$ws.Range("K76").Value = 3632.8333
$ws.Range("M76").Value = -3317.8333
$ws.Range("H77").Value = 5237.25
$ws.Range("I77").Value = 5237.25
$ws.Range("K77").Value = 26186.25
$ws.Range("M77").Value = -21506.25
$ws.Range("H79").Value = 3828.1428
$ws.Range("I79").Value = 3632.8333
$ws.Range("K79").Value = 3632.8333
$ws.Range("M79").Value = -2540.8333
$ws.Range("H82").Value = 3348.6667
$ws.Range("I82").Value = 2999.5
$ws.Range("J82").Value = 4047
$ws.Range("K82").Value = 8998.5
$ws.Range("L82").Value = 12141
$ws.Range("M82").Value = -8592.5
$ws.Range("N82").Value = -12953
$ws.Range("H85").Value = 3348.6667
$ws.Range("I85").Value = 2999.5
$ws.Range("J85").Value = 4047
$ws.Range("K85").Value = 8998.5
$ws.Range("L85").Value = 12141
$ws.Range("M85").Value = -7594.5
$ws.Range("N85").Value = -14949
$ws.Range("H86").Value = 4280.8125
$ws.Range("I86").Value = 2513.3333
$ws.Range("J86").Value = 5341.3
$ws.Range("K86").Value = 2513.3333
$ws.Range("L86").Value = 5341.3
$ws.Range("M86").Value = -1390.3333
$ws.Range("N86").Value = -7587.3
$ws.Range("H89").Value = 4280.8125
$ws.Range("I89").Value = 2513.3333
$ws.Range("J89").Value = 5341.3
$ws.Range("K89").Value = 12566.6665
$ws.Range("L89").Value = 26706.5
$ws.Range("M89").Value = -6950.666499999999
$ws.Range("N89").Value = -37938.5
$ws.Range("H97").Value = 5142.4287
$ws.Range("J97").Value = 5142.4287
$ws.Range("L97").Value = 15427.2861
$ws.Range("N97").Value = -16419.2861
$ws.Range("H100").Value = 3952.5
$ws.Range("I100").Value = 2905
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 2905
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -2364
$ws.Range("N100").Value = -6082
$ws.Range("H112").Value = 3056.875
$ws.Range("I112").Value = 1398
$ws.Range("J112").Value = 3293.8572
$ws.Range("K112").Value = 4194
$ws.Range("L112").Value = 9881.571599999999
$ws.Range("M112").Value = -3086
$ws.Range("N112").Value = -12097.5716
$ws.Range("H132").Value = 2052.8948
$ws.Range("I132").Value = 1944.7222
$ws.Range("K132").Value = 5834.1666
$ws.Range("M132").Value = -3304.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 7797.5
$ws.Range("I28").Value = 7797.5
$ws.Range("K28").Value = 7797.5
$ws.Range("M28").Value = -7605.5
$ws.Range("H37").Value = 8999.375
$ws.Range("I37").Value = 6583.3335
$ws.Range("J37").Value = 16247.5
$ws.Range("K37").Value = 6583.3335
$ws.Range("L37").Value = 16247.5
$ws.Range("M37").Value = -6310.3335
$ws.Range("N37").Value = -16793.5
$ws.Range("H99").Value = 7797.5
$ws.Range("I99").Value = 7797.5
$ws.Range("K99").Value = 7797.5
$ws.Range("M99").Value = -4802.5
$ws.Range("H122").Value = 1997.5
$ws.Range("I122").Value = 1997.5
$ws.Range("K122").Value = 5992.5
$ws.Range("M122").Value = -3542.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 11432
$ws.Range("I97").Value = 10242.667
$ws.Range("K97").Value = 10242.667
$ws.Range("M97").Value = -9251.666999999999
$ws.Range("H99").Value = 3573
$ws.Range("I99").Value = 3573
$ws.Range("K99").Value = 3573
$ws.Range("M99").Value = -2075

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("N92").Value = 0
$ws.Range("H122").Value = 1053.5555
$ws.Range("I122").Value = 1060.25
$ws.Range("K122").Value = 3180.75
$ws.Range("M122").Value = -730.75
$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -54910
$ws.Range("H132").Value = 1632.8064
$ws.Range("I132").Value = 1607.4828
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 4822.4484
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2292.4484
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 173.06667
$ws.Range("I23").Value = 80.22221999999999
$ws.Range("K23").Value = 240.66666
$ws.Range("M23").Value = -5.666659999999979

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3000.2
$ws.Range("J68").Value = 2999.5
$ws.Range("L68").Value = 2999.5
$ws.Range("N68").Value = -4497.5
$ws.Range("H71").Value = 3000.2
$ws.Range("J71").Value = 2999.5
$ws.Range("L71").Value = 14997.5
$ws.Range("N71").Value = -22485.5
$ws.Range("H82").Value = 538
$ws.Range("I82").Value = 450.66666
$ws.Range("K82").Value = 450.66666
$ws.Range("M82").Value = -89.66665999999998
$ws.Range("H85").Value = 538
$ws.Range("I85").Value = 450.66666
$ws.Range("K85").Value = 450.66666
$ws.Range("M85").Value = 797.33334
$ws.Range("H93").Value = 568
$ws.Range("I93").Value = 568
$ws.Range("K93").Value = 568
$ws.Range("M93").Value = 680
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H99").Value = 17078
$ws.Range("I99").Value = 17078
$ws.Range("K99").Value = 17078
$ws.Range("M99").Value = -14083

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2574.75
$ws.Range("J62").Value = 2574.75
$ws.Range("L62").Value = 2574.75
$ws.Range("N62").Value = -3822.75
$ws.Range("H65").Value = 2574.75
$ws.Range("J65").Value = 2574.75
$ws.Range("L65").Value = 12873.75
$ws.Range("N65").Value = -19113.75
$ws.Range("H132").Value = 902.625
$ws.Range("I132").Value = 954.3333
$ws.Range("K132").Value = 2862.9999
$ws.Range("M132").Value = -332.9998999999998
